$wb = $excel.ActiveWorkbook

$changes = @(
    @("ALC", "H98", 2816.9092),
    @("ALC", "I98", 2599.3),
    @("ALC", "K98", 2599.3),
    @("ALC", "M98", -1101.3),
    @("ALC", "H113", 3120.0625),
    @("ALC", "I113", 3054.5454),
    @("ALC", "K113", 3054.5454),
    @("ALC", "M113", 199.4546),
    @("ALC", "H122", 2816.9092),
    @("ALC", "I122", 2599.3),
    @("ALC", "K122", 7797.900000000001),
    @("ALC", "M122", -5347.900000000001),
    @("ALC", "H125", 809.53845),
    @("ALC", "I125", 804.0833),
    @("ALC", "K125", 7236.7497),
    @("ALC", "M125", -4776.7497),
    @("ALC", "H131", 6655.9),
    @("ALC", "I131", 1639.75),
    @("ALC", "K131", 4919.25),
    @("ALC", "M131", 120.75),
    @("ALC", "H138", 14406),
    @("ALC", "J138", 15340.143),
    @("ALC", "L138", 46020.429),
    @("ALC", "N138", -56300.429),
    @("ARM", "H7", 34999),
    @("ARM", "J7", 34999),
    @("ARM", "L7", 34999),
    @("ARM", "N7", -35227),
    @("ARM", "H8", 757501.5),
    @("ARM", "I8", 3000000),
    @("ARM", "K8", 3000000),
    @("ARM", "M8", -2999856),
    @("ARM", "H10", 0),
    @("ARM", "I10", 0),
    @("ARM", "K10", 0),
    @("ARM", "M10", $null),
    @("ARM", "H12", 2115.25),
    @("ARM", "J12", 3494.5),
    @("ARM", "L12", 3494.5),
    @("ARM", "N12", -3840.5),
    @("ARM", "H16", 3639.2),
    @("ARM", "I16", 799),
    @("ARM", "K16", 799),
    @("ARM", "M16", -512),
    @("ARM", "H30", 4249.25),
    @("ARM", "I30", 4249.25),
    @("ARM", "K30", 4249.25),
    @("ARM", "M30", -4099.25),
    @("ARM", "H32", 15201.581),
    @("ARM", "I32", 6234.0835),
    @("ARM", "J32", 27618.115),
    @("ARM", "K32", 6234.0835),
    @("ARM", "L32", 27618.115),
    @("ARM", "M32", -5947.0835),
    @("ARM", "N32", -28192.115),
    @("ARM", "H45", 1970.3334),
    @("ARM", "H54", 32500),
    @("ARM", "J54", 32500),
    @("ARM", "L54", 32500),
    @("ARM", "N54", -34038),
    @("ARM", "H61", 1451.6875),
    @("ARM", "I61", 1409.0714),
    @("ARM", "J61", 1750),
    @("ARM", "K61", 1409.0714),
    @("ARM", "L61", 1750),
    @("ARM", "M61", -1197.0714),
    @("ARM", "N61", -2174),
    @("ARM", "H122", 315024.72),
    @("ARM", "I122", 436522.47),
    @("ARM", "K122", 1309567.41),
    @("ARM", "M122", -1307117.41),
    @("ARM", "H136", 1451.6875),
    @("ARM", "I136", 1409.0714),
    @("ARM", "J136", 1750),
    @("ARM", "K136", 4227.2142),
    @("ARM", "L136", 5250),
    @("ARM", "M136", -1677.2142),
    @("ARM", "N136", -10350),
    @("BSM", "H20", 3629.4546),
    @("BSM", "J20", 7871.6665),
    @("BSM", "L20", 7871.6665),
    @("BSM", "N20", -8365.666499999999),
    @("BSM", "H107", 1865.4615),
    @("BSM", "I107", 1386.5454),
    @("BSM", "J107", 4499.5),
    @("BSM", "K107", 1386.5454),
    @("BSM", "L107", 4499.5),
    @("BSM", "M107", 533.4546),
    @("BSM", "N107", -8339.5),
    @("CRP", "H31", 4823.7646),
    @("CRP", "I31", 2841.5),
    @("CRP", "J31", 5433.6924),
    @("CRP", "K31", 2841.5),
    @("CRP", "L31", 5433.6924),
    @("CRP", "M31", -2546.5),
    @("CRP", "N31", -6023.6924),
    @("CRP", "H34", 4823.7646),
    @("CRP", "I34", 2841.5),
    @("CRP", "J34", 5433.6924),
    @("CRP", "K34", 2841.5),
    @("CRP", "L34", 5433.6924),
    @("CRP", "M34", -2639.5),
    @("CRP", "N34", -5837.6924),
    @("CRP", "H86", 15652.6),
    @("CRP", "I86", 9750),
    @("CRP", "K86", 9750),
    @("CRP", "M86", -8627),
    @("CRP", "H89", 15652.6),
    @("CRP", "I89", 9750),
    @("CRP", "K89", 48750),
    @("CRP", "M89", -43134),
    @("CRP", "H118", 200742),
    @("CRP", "J118", 200742),
    @("CRP", "L118", 200742),
    @("CRP", "N118", -204056),
    @("CUL", "H122", 1168.875),
    @("CUL", "I122", 1058.6666),
    @("CUL", "K122", 9527.999400000001),
    @("CUL", "M122", -7077.999400000001),
    @("CUL", "H131", 1670.8),
    @("CUL", "J131", 2048.9092),
    @("CUL", "L131", 6146.7276),
    @("CUL", "N131", -16226.7276),
    @("CUL", "H132", 4285.5713),
    @("CUL", "J132", 8333.333000000001),
    @("CUL", "L132", 74999.997),
    @("CUL", "N132", -80059.997),
    @("CUL", "H137", 3766.25),
    @("CUL", "I137", 2875.7144),
    @("CUL", "J137", 10000),
    @("CUL", "K137", 8627.143199999999),
    @("CUL", "L137", 30000),
    @("CUL", "M137", -3527.143199999999),
    @("CUL", "N137", -40200),
    @("CUL", "H139", 6486.615),
    @("CUL", "J139", 9332.666999999999),
    @("CUL", "L139", 27998.001),
    @("CUL", "N139", -38278.001),
    @("GSM", "H70", 9964.666999999999),
    @("GSM", "J70", 11499.5),
    @("GSM", "L70", 11499.5),
    @("GSM", "N70", -12039.5),
    @("GSM", "H73", 9964.666999999999),
    @("GSM", "J73", 11499.5),
    @("GSM", "L73", 11499.5),
    @("GSM", "N73", -13371.5),
    @("GSM", "H122", 649840.5),
    @("GSM", "J122", 1432141.1),
    @("GSM", "L122", 4296423.300000001),
    @("GSM", "N122", -4301323.300000001),
    @("LTW", "H16", 6438.643),
    @("LTW", "I16", 6438.643),
    @("LTW", "K16", 6438.643),
    @("LTW", "M16", -6268.643),
    @("LTW", "H55", 695.5417),
    @("LTW", "I55", 552.5789),
    @("LTW", "K55", 552.5789),
    @("LTW", "M55", -379.5789),
    @("LTW", "H82", 1922.5883),
    @("LTW", "I82", 2107.7),
    @("LTW", "K82", 2107.7),
    @("LTW", "M82", -1746.7),
    @("LTW", "H85", 1922.5883),
    @("LTW", "I85", 2107.7),
    @("LTW", "K85", 2107.7),
    @("LTW", "M85", -859.6999999999998),
    @("LTW", "H132", 4747.5415),
    @("LTW", "I132", 3116.8),
    @("LTW", "K132", 9350.400000000001),
    @("LTW", "M132", -6820.400000000001),
    @("LTW", "H136", 4151.1665),
    @("LTW", "I136", 3981.4),
    @("LTW", "K136", 11944.2),
    @("LTW", "M136", -9394.200000000001),
    @("WVR", "H6", 1481.25),
    @("WVR", "J6", 308.33334),
    @("WVR", "L6", 308.33334),
    @("WVR", "N6", -538.33334),
    @("WVR", "H22", 30000),
    @("WVR", "J22", 30000),
    @("WVR", "L22", 30000),
    @("WVR", "N22", -30586),
    @("WVR", "H107", 1492.7778),
    @("WVR", "I107", 776.4286),
    @("WVR", "K107", 2329.2858),
    @("WVR", "M107", -409.2857999999997),
    @("WVR", "H122", 1427.75),
    @("WVR", "I122", 1427.9412),
    @("WVR", "K122", 4283.8236),
    @("WVR", "M122", -1833.8236),
    @("WVR", "H132", 1535.3334),
    @("WVR", "I132", 1139.1818),
    @("WVR", "J132", 2624.75),
    @("WVR", "K132", 3417.5454),
    @("WVR", "L132", 7874.25),
    @("WVR", "M132", -887.5454),
    @("WVR", "N132", -12934.25)
)

foreach ($chg in $changes) {
    $ws = $wb.Worksheets.Item($chg[0])
    $cellRef = $chg[1]
    $val = $chg[2]
    if ($null -eq $val) {
        $ws.Range($cellRef).ClearContents()
    } else {
        $ws.Range($cellRef).Value = $val
    }
}

Write-Host "Applied $($changes.Count) cell updates"